# "Added BoB to presentations"
# Adds a new "Bayes on the Beach" presentation entry to the top of the
# presentations table, and updates the existing "Contributed session
# (planned)" entry (for the Australian Statistical Conference) to drop the
# "(planned)" qualifier now that it has taken place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("presentations")

# Insert a new row above the current row 2 (first data row), pushing all
# existing presentation rows down by one.
$ws.Rows.Item(2).Insert()

# The entry that used to be row 2 (Australian Statistical Conference talk)
# is now row 3; the conference has now occurred, so drop "(planned)".
$ws.Cells.Item(3, 2).Value = "Contributed session"

# Populate the newly inserted row with the new "Bayes on the Beach" entry.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Poster  presentation"
$ws.Cells.Item(2, 4).Value = "Bayes on the Beach"
$ws.Cells.Item(2, 5).Value = " February 2024"
$ws.Cells.Item(2, 6).Value = "Gold Coast, Australia"
$ws.Cells.Item(2, 3).Value = "Generalising the Shared Component Model - The Health Determinants for Cancer Indices for Areas"

# Update the sort range to include the newly added row.
$ws.Sort.SortFields.Clear()
$sortRange = $ws.Range("B6:H17")
$keyRange = $ws.Range("H6:H17")
$ws.Sort.SortFields.Add($keyRange, 0, 2, 0, 1) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Make "presentations" the active sheet/tab, matching the workbook view
# change recorded for this edit.
$ws.Select()
$ws.Range("H9").Select()
